$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swaps (rows display swapped order due to updated ranking) ---
# China <-> Panama
$ws.Cells.Item(37, 1).Value = "Panama"
$ws.Cells.Item(38, 1).Value = "China"

# Guyana <-> Trinidad yTobago
$ws.Cells.Item(162, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(163, 1).Value = "Guyana"

# --- Updated timestamp ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 23 de Agosto de 2020 a las 01:32"

# --- Updated numeric statistics per country row ---
$ws.Cells.Item(4, 2).Value = 5838632
$ws.Cells.Item(4, 3).Value = 41895
$ws.Cells.Item(4, 4).Value = 3144164
$ws.Cells.Item(4, 5).Value = 2514328
$ws.Cells.Item(4, 7).Value = 940
$ws.Cells.Item(4, 8).Value = 180140
$ws.Cells.Item(5, 2).Value = 3582698
$ws.Cells.Item(5, 3).Value = 46210
$ws.Cells.Item(5, 5).Value = 797666
$ws.Cells.Item(5, 7).Value = 823
$ws.Cells.Item(5, 8).Value = 114277
$ws.Cells.Item(11, 2).Value = 533103
$ws.Cells.Item(11, 3).Value = 10965
$ws.Cells.Item(11, 4).Value = 359792
$ws.Cells.Item(11, 5).Value = 156343
$ws.Cells.Item(11, 7).Value = 400
$ws.Cells.Item(11, 8).Value = 16968
$ws.Cells.Item(15, 2).Value = 336802
$ws.Cells.Item(15, 3).Value = 7759
$ws.Cells.Item(15, 5).Value = 84173
$ws.Cells.Item(15, 7).Value = 118
$ws.Cells.Item(15, 8).Value = 6848
$ws.Cells.Item(23, 2).Value = 233857
$ws.Cells.Item(23, 3).Value = 836
$ws.Cells.Item(23, 5).Value = 15576
$ws.Cells.Item(27, 2).Value = 124629
$ws.Cells.Item(27, 3).Value = 257
$ws.Cells.Item(27, 4).Value = 110878
$ws.Cells.Item(27, 5).Value = 4680
$ws.Cells.Item(37, 2).Value = 85480
$ws.Cells.Item(37, 3).Value = 1088
$ws.Cells.Item(37, 4).Value = 60528
$ws.Cells.Item(37, 5).Value = 23074
$ws.Cells.Item(37, 7).Value = 19
$ws.Cells.Item(37, 8).Value = 1878
$ws.Cells.Item(38, 2).Value = 84939
$ws.Cells.Item(38, 3).Value = 22
$ws.Cells.Item(38, 4).Value = 79851
$ws.Cells.Item(38, 5).Value = 454
$ws.Cells.Item(38, 8).Value = 4634
$ws.Cells.Item(52, 2).Value = 51905
$ws.Cells.Item(52, 3).Value = 601
$ws.Cells.Item(52, 4).Value = 38767
$ws.Cells.Item(52, 5).Value = 12141
$ws.Cells.Item(52, 7).Value = 1
$ws.Cells.Item(52, 8).Value = 997
$ws.Cells.Item(74, 2).Value = 21790
$ws.Cells.Item(74, 3).Value = 239
$ws.Cells.Item(74, 5).Value = 5292
$ws.Cells.Item(85, 2).Value = 12682
$ws.Cells.Item(85, 3).Value = 59
$ws.Cells.Item(85, 4).Value = 6492
$ws.Cells.Item(85, 5).Value = 5375
$ws.Cells.Item(85, 7).Value = 3
$ws.Cells.Item(85, 8).Value = 815
$ws.Cells.Item(89, 2).Value = 10299
$ws.Cells.Item(89, 3).Value = 24
$ws.Cells.Item(89, 5).Value = 885
$ws.Cells.Item(102, 2).Value = 7762
$ws.Cells.Item(102, 3).Value = 58
$ws.Cells.Item(102, 5).Value = 669
$ws.Cells.Item(104, 2).Value = 6660
$ws.Cells.Item(104, 3).Value = 96
$ws.Cells.Item(104, 4).Value = 4113
$ws.Cells.Item(104, 5).Value = 2521
$ws.Cells.Item(110, 4).Value = 3795
$ws.Cells.Item(110, 5).Value = 1048
$ws.Cells.Item(114, 2).Value = 4311
$ws.Cells.Item(114, 3).Value = 34
$ws.Cells.Item(114, 4).Value = 3333
$ws.Cells.Item(114, 5).Value = 894
$ws.Cells.Item(114, 7).Value = 2
$ws.Cells.Item(114, 8).Value = 84
$ws.Cells.Item(118, 2).Value = 3569
$ws.Cells.Item(118, 3).Value = 109
$ws.Cells.Item(118, 4).Value = 2559
$ws.Cells.Item(118, 5).Value = 953
$ws.Cells.Item(118, 7).Value = 1
$ws.Cells.Item(118, 8).Value = 57
$ws.Cells.Item(119, 2).Value = 3455
$ws.Cells.Item(119, 3).Value = 43
$ws.Cells.Item(119, 4).Value = 2538
$ws.Cells.Item(119, 5).Value = 880
$ws.Cells.Item(138, 2).Value = 2115
$ws.Cells.Item(138, 3).Value = 20
$ws.Cells.Item(138, 5).Value = 371
$ws.Cells.Item(147, 2).Value = 1521
$ws.Cells.Item(147, 3).Value = 5
$ws.Cells.Item(147, 4).Value = 1264
$ws.Cells.Item(147, 5).Value = 215
$ws.Cells.Item(162, 2).Value = 930
$ws.Cells.Item(162, 3).Value = 66
$ws.Cells.Item(162, 4).Value = 165
$ws.Cells.Item(162, 5).Value = 752
$ws.Cells.Item(162, 8).Value = 13
$ws.Cells.Item(163, 2).Value = 925
$ws.Cells.Item(163, 3).Value = 44
$ws.Cells.Item(163, 4).Value = 433
$ws.Cells.Item(163, 5).Value = 461
$ws.Cells.Item(163, 8).Value = 31
